$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    4  = -7
    6  = -8
    7  = -7
    8  = -14
    9  = -9
    10 = -8
    11 = -1
    13 = -5
    16 = -2
    17 = -2
    18 = -4
    19 = -3
    21 = -2
    22 = -1
    23 = 7
    24 = -4
    25 = 1
    28 = -4
    29 = -3
    33 = -3
    34 = 2
    35 = -9
    36 = -3
    37 = -4
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
